$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (plain decimal-looking strings), so
# they round-trip as text just like the rest of column D.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '60.831.67'
$ws.Range('E2').Value = '  -3.29%  '
$ws.Range('D3').Value = '2.906.42'
$ws.Range('E3').Value = '  -3.95%  '
$ws.Range('D5').Value = '590.04'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '144.75'
$ws.Range('E6').Value = '  -5.30%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = '2.904.56'
$ws.Range('E9').Value = '  -3.92%  '
$ws.Range('E10').Value = '  -4.57%  '
$ws.Range('E11').Value = '  -4.13%  '
$ws.Range('E12').Value = '  -4.05%  '
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('D14').Value = '33.44'
$ws.Range('E14').Value = '  -6.00%  '
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').Value = '3.388.80'
$ws.Range('D17').Value = '60.774.86'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('E18').Value = '  -4.95%  '
$ws.Range('D19').Value = '2.906.21'
$ws.Range('E19').Value = '  -3.94%  '
$ws.Range('D20').Value = '429.71'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('E21').Value = '  -4.79%  '
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('E23').Value = '  -5.69%  '
$ws.Range('D24').Value = '81.89'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').Value = '10.79'
$ws.Range('E25').Value = '  -5.95%  '
$ws.Range('E26').Value = '  -4.28%  '
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('E32').Value = '  -6.31%  '
$ws.Range('D33').Value = '26.53'
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').Value = '0.0₃0854'
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('E36').Value = '  -3.35%  '
$ws.Range('E37').Value = '  -4.61%  '
$ws.Range('D38').Value = '3.01'
$ws.Range('E38').Value = '  -4.76%  '
$ws.Range('D39').Value = '49.57'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('E41').Value = '  -4.89%  '
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('D43').Value = '0.292'
$ws.Range('E43').Value = '  -4.62%  '
$ws.Range('E44').Value = '  -10.05%  '
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('D46').Value = '373.11'
$ws.Range('E46').Value = '  -4.88%  '
$ws.Range('D47').Value = '2.700.09'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').Value = '131.28'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D50').Value = '24.03'
$ws.Range('E50').Value = '  -9.33%  '
$ws.Range('E51').Value = '  -2.10%  '
